$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.639.63"
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.088.09"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "345.03"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5162"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4391"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09225"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.97"
$ws.Range("E10").Value = "  -3.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.176"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.36"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.088.75"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.213"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.735"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.01"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001157"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.04"
$ws.Range("E19").Value = "  +8.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06625"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.195"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.669.45"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.68"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.332.93"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.88"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.34"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.77"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.145"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.627"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.189"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.969"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.104"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.29"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02571"
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06725"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2274"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.56"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.289"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6641"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.19"
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.331"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.624"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.65"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000331"
$ws.Range("E50").Value = "  -8.33%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.167"
$ws.Range("E51").Value = "  -2.29%  "
